$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as text, avoiding Excel auto-converting numeric-looking strings to numbers
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.031.64"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.635.95"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.52%  "
Set-TextValue $ws.Range("D5") "214.75"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.50%  "
Set-TextValue $ws.Range("D8") "0.251"
$ws.Range("E8").Value = "  -1.89%  "
Set-TextValue $ws.Range("D9") "0.0627"
$ws.Range("E9").Value = "  -1.48%  "
Set-TextValue $ws.Range("D10") "18.69"
$ws.Range("E10").Value = "  -4.37%  "
Set-TextValue $ws.Range("D11") "0.0794"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.704.96"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("D13").Value = "1.865.72"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  -1.49%  "
Set-TextValue $ws.Range("D15") "0.532"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "0.0₃0748"
$ws.Range("E16").Value = "  -1.96%  "
Set-TextValue $ws.Range("D17") "62.21"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "26.048.11"
$ws.Range("E18").Value = "  +0.41%  "
Set-TextValue $ws.Range("D20") "191.04"
Set-TextValue $ws.Range("D21") "4.27"
$ws.Range("E21").Value = "  -1.93%  "
Set-TextValue $ws.Range("D22") "9.61"
$ws.Range("E22").Value = "  -3.02%  "
Set-TextValue $ws.Range("D23") "6.15"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("E24").Value = "  -0.03%  "
Set-TextValue $ws.Range("D25") "143.59"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D26") "1.01"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D27") "1.78"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -0.67%  "
Set-TextValue $ws.Range("D36") "0.879"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").Value = "1.130.69"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -2.92%  "
Set-TextValue $ws.Range("D40") "0.0155"
$ws.Range("E40").Value = "  -0.85%  "
Set-TextValue $ws.Range("D41") "98.94"
$ws.Range("E41").Value = "  -0.33%  "
Set-TextValue $ws.Range("D42") "0.791"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("E44").Value = "  -1.03%  "
Set-TextValue $ws.Range("D45") "55.56"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  +0.00%  "
Set-TextValue $ws.Range("D49") "7.59"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  +0.38%  "
Set-TextValue $ws.Range("D51") "0.0929"
$ws.Range("E51").Value = "  -3.18%  "
